# Update the dSF column (F) values for specific rows to reflect a data
# repull / recalculation, per the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    4  = 1
    11 = 0
    15 = 6
    20 = 2
    23 = -4
    27 = -3
    33 = 2
    34 = 2
    38 = 0
    44 = 2
    45 = -9
    47 = -9
    49 = -4
    50 = -1
    57 = -2
    62 = 3
    65 = -4
    68 = 0
    69 = 3
    70 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
